# Update match-odds values that changed in the 2024-10-11 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.35
$ws.Range("H2").Value = 2.7

# Row 3
$ws.Range("G3").Value = 2.75
$ws.Range("H3").Value = 2.8
$ws.Range("I3").Value = 2.9

# Row 4
$ws.Range("M4").Value = 1.11
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.38

# Row 5
$ws.Range("G5").Value = 1.62
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 5.5
$ws.Range("J5").Value = 2.2
$ws.Range("L5").Value = 5.5
$ws.Range("M5").Value = 1.05
$ws.Range("O5").Value = 1.29
$ws.Range("Q5").Value = 1.98
$ws.Range("R5").Value = 1.88
$ws.Range("X5").Value = 7.5
$ws.Range("Z5").Value = 12
$ws.Range("AB5").Value = 26
$ws.Range("AE5").Value = 17
$ws.Range("AG5").Value = 13
$ws.Range("AH5").Value = 26
$ws.Range("AI5").Value = 17
$ws.Range("AJ5").Value = 51
$ws.Range("AL5").Value = 41
$ws.Range("AM5").Value = 301
$ws.Range("AN5").Value = 3.6
$ws.Range("AO5").Value = 8.5
$ws.Range("AY5").Value = 34
$ws.Range("AZ5").Value = 101
$ws.Range("BA5").Value = 126
$ws.Range("BB5").Value = 251

# Row 6
$ws.Range("G6").Value = 2.55
$ws.Range("I6").Value = 2.7
$ws.Range("J6").Value = 3.1
$ws.Range("L6").Value = 3.25
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11
$ws.Range("O6").Value = 1.25
$ws.Range("Q6").Value = 1.85
$ws.Range("R6").Value = 2
$ws.Range("U6").Value = 1.67
$ws.Range("V6").Value = 2.1
$ws.Range("W6").Value = 9.5
$ws.Range("AA6").Value = 19
$ws.Range("AC6").Value = 11
$ws.Range("AH6").Value = 15
$ws.Range("AJ6").Value = 29
$ws.Range("AL6").Value = 29
$ws.Range("AO6").Value = 13
$ws.Range("AY6").Value = 23
$ws.Range("AZ6").Value = 51
$ws.Range("BA6").Value = 67

# Row 7
$ws.Range("G7").Value = 2.3
$ws.Range("H7").Value = 2.87
$ws.Range("I7").Value = 3.25
$ws.Range("J7").Value = 2.95
$ws.Range("L7").Value = 3.7
$ws.Range("N7").Value = 6.7
$ws.Range("P7").Value = 2.5
$ws.Range("Q7").Value = 2.18
$ws.Range("R7").Value = 1.53
$ws.Range("U7").Value = 1.85
$ws.Range("V7").Value = 1.75
$ws.Range("W7").Value = 6.3
$ws.Range("X7").Value = 10.25
$ws.Range("Y7").Value = 9.25
$ws.Range("Z7").Value = 24
$ws.Range("AA7").Value = 22
$ws.Range("AB7").Value = 37
$ws.Range("AG7").Value = 8.5
$ws.Range("AH7").Value = 17
$ws.Range("AI7").Value = 11.25
$ws.Range("AJ7").Value = 45
$ws.Range("AK7").Value = 30
$ws.Range("AL7").Value = 40
$ws.Range("AN7").Value = 4.1
$ws.Range("AO7").Value = 12.5
$ws.Range("AP7").Value = 21
$ws.Range("AQ7").Value = 55
$ws.Range("AR7").Value = 100
$ws.Range("AS7").Value = 300
$ws.Range("AW7").Value = 5
$ws.Range("AX7").Value = 17.5
$ws.Range("AY7").Value = 24
$ws.Range("AZ7").Value = 90

# Row 10
$ws.Range("H10").Value = 2.7
$ws.Range("K10").Value = 1.91
$ws.Range("O10").Value = 1.53
$ws.Range("P10").Value = 2.38
$ws.Range("Q10").Value = 2.7
$ws.Range("R10").Value = 1.44
$ws.Range("AF10").Value = 67
